$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the new bug-report row off the existing row's formatting so
# every column picks up the right number format / alignment / wrap.
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)

# Column E ("Status") on row 4 is bold (it highlights the current bug's
# status); row 5 should use the plain vertical-center style like the
# other non-wrapped columns, so re-stamp its format from C4.
$ws.Range("C4").Copy()
$ws.Range("E5").PasteSpecial(-4122)

# Populate the new bug report: character continues to slide after death.
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 45338
$ws.Range("C5").Value = "Daniel Lee"
$ws.Range("D5").Value = "P2"
$ws.Range("E5").Value = "open"
$ws.Range("F5").Value = "Not yet assigned"
$ws.Range("G5").Value = "The character is expected to come to a stop upon death with no animations playing. However, currently, the character continues to move horizontally after death, indicating that horizontal velocity is not being reset."
$ws.Range("H5").Value = "When the character dies (for example, by being hit by a falling boulder), all movement and animations should cease."
$ws.Range("I5").Value = "After the character is hit and the death event is triggered, the character continues to slide horizontally."
$ws.Range("J5").Value = "1. Allow the character to be hit by the boulder to trigger the death sequence.`n2. Observe the character's behavior following the death event."

$ws.Rows.Item(5).RowHeight = 90

# Move the active selection/scroll position like the authored workbook.
[void]$ws.Range("I11").Select()
